$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 665.3333
$ws.Cells.Item(4, 9).Value = 665.3333
$ws.Cells.Item(4, 11).Value = 665.3333
$ws.Cells.Item(4, 13).Value = -551.3333

$ws.Cells.Item(5, 8).Value = 173.75
$ws.Cells.Item(5, 9).Value = 173.75
$ws.Cells.Item(5, 11).Value = 173.75
$ws.Cells.Item(5, 13).Value = -58.75

$ws.Cells.Item(18, 8).Value = 663.5
$ws.Cells.Item(18, 9).Value = 663.5
$ws.Cells.Item(18, 11).Value = 663.5
$ws.Cells.Item(18, 13).Value = -379.5

$ws.Cells.Item(86, 8).Value = 3349.6
$ws.Cells.Item(86, 9).Value = 3177
$ws.Cells.Item(86, 11).Value = 3177
$ws.Cells.Item(86, 13).Value = -2054

$ws.Cells.Item(89, 8).Value = 3349.6
$ws.Cells.Item(89, 9).Value = 3177
$ws.Cells.Item(89, 11).Value = 15885
$ws.Cells.Item(89, 13).Value = -10269

$ws.Cells.Item(107, 8).Value = 1565.6061
$ws.Cells.Item(107, 9).Value = 1353
$ws.Cells.Item(107, 10).Value = 2230
$ws.Cells.Item(107, 11).Value = 1353
$ws.Cells.Item(107, 12).Value = 2230
$ws.Cells.Item(107, 13).Value = 567
$ws.Cells.Item(107, 14).Value = -6070

$ws.Cells.Item(111, 8).Value = 1058.3334
$ws.Cells.Item(111, 10).Value = 475
$ws.Cells.Item(111, 12).Value = 1425
$ws.Cells.Item(111, 14).Value = -7559

$ws.Cells.Item(135, 8).Value = 19349.621
$ws.Cells.Item(135, 9).Value = 1804.12
$ws.Cells.Item(135, 11).Value = 16237.08
$ws.Cells.Item(135, 13).Value = -13702.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1070.5
$ws.Cells.Item(5, 9).Value = 1270.6
$ws.Cells.Item(5, 10).Value = 70
$ws.Cells.Item(5, 11).Value = 1270.6
$ws.Cells.Item(5, 12).Value = 70
$ws.Cells.Item(5, 13).Value = -1158.6
$ws.Cells.Item(5, 14).Value = -294

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1070.5
$ws.Cells.Item(4, 9).Value = 1270.6
$ws.Cells.Item(4, 10).Value = 70
$ws.Cells.Item(4, 11).Value = 1270.6
$ws.Cells.Item(4, 12).Value = 70
$ws.Cells.Item(4, 13).Value = -1155.6
$ws.Cells.Item(4, 14).Value = -300

$ws.Cells.Item(94, 8).Value = 1327.826
$ws.Cells.Item(94, 9).Value = 1453.0625
$ws.Cells.Item(94, 11).Value = 1453.0625
$ws.Cells.Item(94, 13).Value = -1002.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 237.93333
$ws.Cells.Item(7, 9).Value = 241.81818
$ws.Cells.Item(7, 10).Value = 227.25
$ws.Cells.Item(7, 11).Value = 241.81818
$ws.Cells.Item(7, 12).Value = 227.25
$ws.Cells.Item(7, 13).Value = -128.81818
$ws.Cells.Item(7, 14).Value = -453.25

$ws.Cells.Item(19, 8).Value = 11565585
$ws.Cells.Item(19, 9).Value = 14133548
$ws.Cells.Item(19, 11).Value = 14133548
$ws.Cells.Item(19, 13).Value = -14133378

$ws.Cells.Item(22, 8).Value = 682.9
$ws.Cells.Item(22, 9).Value = 689.8
$ws.Cells.Item(22, 10).Value = 676
$ws.Cells.Item(22, 11).Value = 689.8
$ws.Cells.Item(22, 12).Value = 676
$ws.Cells.Item(22, 13).Value = -339.8
$ws.Cells.Item(22, 14).Value = -1376

$ws.Cells.Item(24, 8).Value = 11565585
$ws.Cells.Item(24, 9).Value = 14133548
$ws.Cells.Item(24, 11).Value = 14133548
$ws.Cells.Item(24, 13).Value = -14133378

$ws.Cells.Item(26, 8).Value = 30000
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 30000
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 30000
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 14).Value = -30574

$ws.Cells.Item(31, 8).Value = 3337332.2
$ws.Cells.Item(31, 9).Value = 4350846
$ws.Cells.Item(31, 11).Value = 4350846
$ws.Cells.Item(31, 13).Value = -4350551

$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 13).ClearContents()

$ws.Cells.Item(34, 8).Value = 3337332.2
$ws.Cells.Item(34, 9).Value = 4350846
$ws.Cells.Item(34, 11).Value = 4350846
$ws.Cells.Item(34, 13).Value = -4350644

$ws.Cells.Item(92, 8).Value = 20300.5
$ws.Cells.Item(92, 10).Value = 20300.5
$ws.Cells.Item(92, 12).Value = 20300.5
$ws.Cells.Item(92, 14).Value = -25292.5

$ws.Cells.Item(140, 8).Value = 119415.25
$ws.Cells.Item(140, 10).Value = 119415.25
$ws.Cells.Item(140, 12).Value = 119415.25
$ws.Cells.Item(140, 14).Value = -129775.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 548.5
$ws.Cells.Item(8, 9).Value = 548.5
$ws.Cells.Item(8, 11).Value = 1645.5
$ws.Cells.Item(8, 13).Value = -1506.5

$ws.Cells.Item(105, 8).Value = 6767.636
$ws.Cells.Item(105, 9).Value = 5632.6665
$ws.Cells.Item(105, 11).Value = 16897.9995
$ws.Cells.Item(105, 13).Value = -14276.9995

$ws.Cells.Item(121, 8).Value = 114153.5
$ws.Cells.Item(121, 9).Value = 186088.17
$ws.Cells.Item(121, 10).Value = 6251.5
$ws.Cells.Item(121, 11).Value = 558264.51
$ws.Cells.Item(121, 12).Value = 18754.5
$ws.Cells.Item(121, 13).Value = -556954.51
$ws.Cells.Item(121, 14).Value = -21374.5

$ws.Cells.Item(131, 8).Value = 427351.2
$ws.Cells.Item(131, 9).Value = 851396.8
$ws.Cells.Item(131, 10).Value = 3305.6
$ws.Cells.Item(131, 11).Value = 2554190.4
$ws.Cells.Item(131, 12).Value = 9916.799999999999
$ws.Cells.Item(131, 13).Value = -2549150.4
$ws.Cells.Item(131, 14).Value = -19996.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 23069
$ws.Cells.Item(40, 10).Value = 23069
$ws.Cells.Item(40, 12).Value = 23069
$ws.Cells.Item(40, 14).Value = -23371

$ws.Cells.Item(70, 8).Value = 6293.85
$ws.Cells.Item(70, 9).Value = 6631.3335
$ws.Cells.Item(70, 11).Value = 6631.3335
$ws.Cells.Item(70, 13).Value = -6361.3335

$ws.Cells.Item(73, 8).Value = 6293.85
$ws.Cells.Item(73, 9).Value = 6631.3335
$ws.Cells.Item(73, 11).Value = 6631.3335
$ws.Cells.Item(73, 13).Value = -5695.3335

$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 80000
$ws.Cells.Item(25, 9).Value = 80000
$ws.Cells.Item(25, 10).Value = 80000
$ws.Cells.Item(25, 11).Value = 80000
$ws.Cells.Item(25, 12).Value = 80000
$ws.Cells.Item(25, 13).Value = -79770
$ws.Cells.Item(25, 14).Value = -80460

$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 34400
$ws.Cells.Item(42, 10).Value = 34400
$ws.Cells.Item(42, 12).Value = 34400
$ws.Cells.Item(42, 14).Value = -35156

$ws.Cells.Item(140, 8).Value = 95000
$ws.Cells.Item(140, 10).Value = 95000
$ws.Cells.Item(140, 12).Value = 95000
$ws.Cells.Item(140, 14).Value = -105360

$ws.Cells.Item(141, 8).Value = 94993.25
$ws.Cells.Item(141, 10).Value = 94993.25
$ws.Cells.Item(141, 12).Value = 94993.25
$ws.Cells.Item(141, 14).Value = -105353.25
